$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in column C for every data row
# (rows 2-42) from 2025-02-18 (45706) to 2025-02-19 (45707).
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45707
}

# Row 39 and row 40 had their "Beteckning" (A) and "Area (ha)" (G) values
# swapped: A39/G39 <-> A40/G40.
$a39 = $ws.Cells.Item(39, 1).Value2
$g39 = $ws.Cells.Item(39, 7).Value2
$a40 = $ws.Cells.Item(40, 1).Value2
$g40 = $ws.Cells.Item(40, 7).Value2

$ws.Cells.Item(39, 1).Value2 = $a40
$ws.Cells.Item(39, 7).Value2 = $g40
$ws.Cells.Item(40, 1).Value2 = $a39
$ws.Cells.Item(40, 7).Value2 = $g39
